# Procédure de tests système - add audio/software/wifi test rows (S4, S5, S6),
# add two "notes" remarks, rename the "Liste des tests" section header, and
# normalize row 26 formatting. Matches the "Procédure de test Fini" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Section title above the table: "Schéma Électrique" -> "Liste des tests"
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "Liste des tests"

# ---------------------------------------------------------------------------
# 2) Fill in the two "notes" (column H) that were left blank on existing rows
# ---------------------------------------------------------------------------
$ws.Range("H16").Value = "Attention le test pourrais endommagé la carte "
$ws.Range("H22").Value = "Afficher les trois couleur sur chaque pixel"

# ---------------------------------------------------------------------------
# 3) Normalize row 26 formatting (A26/B26/C26 used stray "applyFill" style
#    variants; bring them back in line with the rest of the table - same
#    borders/wrap, no fill override). Values are untouched.
# ---------------------------------------------------------------------------
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A26").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B24").Copy() | Out-Null
$ws.Range("B26:C26").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Three new test rows: S4 (Audio), S5 (Logiciel/panneau), S6 (Wifi)
#    Formats are copied from existing rows first so the new cells reuse the
#    workbook's existing styles, then the values are written.
# ---------------------------------------------------------------------------

# Row 27 - S4 / Test Audio
$ws.Range("A22").Copy() | Out-Null
$ws.Range("A27").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B22").Copy() | Out-Null
$ws.Range("B27").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C21").Copy() | Out-Null
$ws.Range("C27").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D22").Copy() | Out-Null
$ws.Range("D27").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E22").Copy() | Out-Null
$ws.Range("E27").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H16").Copy() | Out-Null
$ws.Range("H27").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A27").Value = "S4"
$ws.Range("B27").Value = "Test Audio"
$ws.Range("C27").Value = "électronique et logiciel"
$ws.Range("D27").Value = "S'assurer que l'enregistrement et le speaker fonctionne"
$ws.Range("E27").Value = "Fonctionnel"
$ws.Rows.Item(27).RowHeight = 28.5

# Row 28 - S5 / options du panneau
$ws.Range("F24").Copy() | Out-Null
$ws.Range("A28:D28").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H16").Copy() | Out-Null
$ws.Range("H28").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A28").Value = "S5"
$ws.Range("B28").Value = "Tester les options du panneau"
$ws.Range("C28").Value = "Logiciel"
$ws.Range("D28").Value = "S'assurer que chaque option du panneau soit fonctionnel"
$ws.Range("E28").Value = 0.9
$ws.Range("E28").NumberFormat = "0%"
$ws.Range("H28").Value = "Option de désactivation des capteurs, changement de code etc."

# Row 29 - S6 / communication Wifi
$ws.Range("F24").Copy() | Out-Null
$ws.Range("A29:E29").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H16").Copy() | Out-Null
$ws.Range("H29").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A29").Value = "S6"
$ws.Range("B29").Value = "Tester la communication Wifi"
$ws.Range("C29").Value = "Logiciel"
$ws.Range("D29").Value = "Communication UDP fonctionnel avec une page web basique"
$ws.Range("E29").Value = "Fonctionnel"
$ws.Range("H29").Value = "Avoir le serveur web fonctionnel pour se connecter avec des appareils"
$ws.Rows.Item(29).RowHeight = 28.5

# ---------------------------------------------------------------------------
# 5) Restore the view: scrolled down with D25 selected
# ---------------------------------------------------------------------------
$ws.Range("D25").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 2
